$data = @(
    @(1.02,1.069596262564414,1.072208229385172,1.082310250929075,1.087425661731959,1.049721509348661,1.07453018593762,1.074903354595115,1.084978800921521,1.090080957938863,1.028219226534989),
    @(1.02,1.071145700132907,1.073459154800193,1.083810287952381,1.088909797004621,1.050138724434446,1.075734016557425,1.075969909923792,1.086295799571925,1.091383061924657,1.028643765535135),
    @(1.02,1.07214677818707,1.074267008536865,1.084779725124916,1.089868912805929,1.050406336410323,1.076511042660754,1.076657902056601,1.087146283050843,1.092223870542934,1.028917271206506),
    @(1.02,1.072567277764295,1.074606258122804,1.085186999559972,1.090271840302838,1.050518279415545,1.076837248002072,1.076946626954748,1.087503425317982,1.092576936465576,1.029031968049303),
    @(1.02,1.072637860995288,1.074663197967509,1.085255366698659,1.090339477106216,1.050537042309651,1.076891992730757,1.07699507556263,1.08756336775441,1.092636193931139,1.02905120950878),
    @(1.02,1.072152398302335,1.074271543061465,1.084785168225129,1.089874297853177,1.050407834399785,1.076515403218195,1.076661761998256,1.087151056770692,1.092228589830177,1.028918804908219),
    @(1.02,1.070120220228058,1.072631315864044,1.082817445627532,1.087927488837407,1.049862998081271,1.074937430126083,1.075264248449352,1.08542424335122,1.090521374208528,1.028362951165024),
    @(1.02,1.066527314215535,1.069728673861391,1.07934064985487,1.084487310360036,1.048884784461355,1.072141754263163,1.072785012466049,1.082368037838703,1.087499426627026,1.027374180162719),
    @(1.02,1.064123489124628,1.067784905997717,1.077015997014792,1.082186923432553,1.048220284457435,1.070267430203518,1.071120656634679,1.080321178486777,1.085475230902013,1.026708616564562),
    @(1.02,1.063080460854707,1.066941093132331,1.076007683942269,1.081189085954198,1.047929581508114,1.069453238142387,1.070397159484684,1.079432544670458,1.084596371474035,1.026418874864942),
    @(1.02,1.062692699667028,1.066627333441672,1.075632883791169,1.080818172372101,1.047821152087778,1.069150412973148,1.070127989970217,1.079102108216591,1.084269559980762,1.026311016541865),
    @(1.02,1.062775891037605,1.066694651011265,1.075713291957368,1.080897747057483,1.047844430954299,1.069215388134867,1.070185747333331,1.079173004303302,1.08433967875515,1.026334163212908),
    @(1.02,1.063048415267209,1.066915164432791,1.075976708359317,1.081158431727391,1.04792062788458,1.069428214683334,1.070374918669275,1.079405238046746,1.084569364605495,1.026409964079435),
    @(1.02,1.06321628187115,1.067050986081334,1.07613897214701,1.081319011811037,1.047967515682351,1.069559291173812,1.070491416129163,1.079548277197657,1.084710833168195,1.026456636225742),
    @(1.02,1.064192665187058,1.067840861195685,1.077082878370861,1.082253108909581,1.048239514607965,1.070321410064543,1.071168612764928,1.080380104393251,1.085533507307008,1.026727812912825),
    @(1.02,1.064804540270699,1.068335748676058,1.077674498228059,1.082838567435495,1.048409334939922,1.07079876626005,1.071592640541492,1.080901257984385,1.086048909092771,1.02689749843683),
    @(1.02,1.065161229305844,1.068624201358257,1.078019414438012,1.083179886770424,1.048508101915399,1.07107694997155,1.071839697132633,1.08120501378331,1.086349306174168,1.026996323933287),
    @(1.02,1.065282816183647,1.068722521407018,1.078136994124438,1.083296239486917,1.04854173042725,1.071171761243733,1.071923891175108,1.08130854886713,1.086451695373943,1.027029995649616),
    @(1.02,1.064738913385276,1.068282673376913,1.077611040214256,1.082775770790596,1.048391144467165,1.07074757638121,1.071547174519781,1.080845366393288,1.085993635006238,1.026879308238235),
    @(1.02,1.062968172977017,1.066850237894862,1.075899146258413,1.08108167417677,1.047898202219163,1.069365553618359,1.070319224363273,1.079336860954268,1.084501737959338,1.02638764913263),
    @(1.02,1.061852900441199,1.065947695825472,1.074821256819835,1.080014948520852,1.047585668166442,1.068494314542606,1.069544669461212,1.078386326657866,1.083561613237493,1.026077160472442),
    @(1.02,1.062444314586896,1.066426334214278,1.07539281667764,1.080580592644376,1.047751596062148,1.068956396098116,1.069955514338535,1.078890422508296,1.084060193818016,1.026241886510371),
    @(1.02,1.064768567990804,1.068306656459271,1.077639714671238,1.082804146412948,1.048399364843711,1.070770707656936,1.071567719497461,1.080870622090677,1.086018611691604,1.026887528067094),
    @(1.02,1.067457635905604,1.070480577914188,1.080240647762946,1.0853778683989,1.049139842441006,1.072866331941258,1.073427961491892,1.083159764098281,1.08828232800347,1.027630915897657)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,3,4,5,6,9,10,11,12,13,14)  # B,C,D,E,F,I,J,K,L,M,N

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($rowNum, $cols[$j]).Value = $rowVals[$j]
    }
}
